$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13, shifting existing rows 13-18 down to 14-19
# (Excel copies formatting, e.g. the date style, from the row above on insert).
$ws.Rows.Item(13).Insert()

# Fill the new row 13 with the new weekly price record.
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 45271
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100112021
$ws.Range("G13").Value = "Ají"
$ws.Range("H13").Value = "Americana (o)"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 10
$ws.Range("K13").Value = 45000
$ws.Range("L13").Value = 45000
$ws.Range("M13").Value = 45000
$ws.Range("N13").Value = "`$/caja 25 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 1800
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
